$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "16th MHz Crystal" -> "ABM7" device name (row 9) ---
$ws.Range("B9").Value = "ABM7"

# --- Re-style the existing quantity-price row 17 (K:N) with the new currency format ---
$ws.Range("K17:N17").NumberFormat = "[$$-409]#,##0.00;[RED]\-[$$-409]#,##0.00"
$ws.Range("K17:N17").HorizontalAlignment = -4108

# --- Add the new BOM row 18 for the 8.06k resistor ---
$ws.Range("B18").Value = "8.06k Resistor"
$ws.Range("C18").Value = "0603 "
$ws.Range("D18").Value = "8.06k"
$ws.Range("E18").Value = "0.1W"
$ws.Range("F18").Value = "RES SMD 8.06K OHM 1% 1/10W 0603"
$ws.Range("G18").Value = "Panasonic electronic Components"
$ws.Range("H18").Value = "ERJ-3EKF8061V"
$ws.Range("I18").Value = "P8.06KHCT-ND"
$ws.Range("K18").Value = 0.1
$ws.Range("L18").Value = "-"
$ws.Range("M18").Value = 0.0114
$ws.Range("N18").Value = 0.00416

# Formats matching the other BOM rows (13, 15, 16)
$ws.Range("B18").NumberFormat = $ws.Range("B17").NumberFormat
$ws.Range("B18").HorizontalAlignment = $ws.Range("B17").HorizontalAlignment

$ws.Range("C18").NumberFormat = $ws.Range("C17").NumberFormat
$ws.Range("C18").HorizontalAlignment = $ws.Range("C17").HorizontalAlignment

$ws.Range("D18").NumberFormat = $ws.Range("D17").NumberFormat
$ws.Range("D18").HorizontalAlignment = $ws.Range("D17").HorizontalAlignment

$ws.Range("E18").NumberFormat = $ws.Range("E17").NumberFormat
$ws.Range("E18").HorizontalAlignment = $ws.Range("E17").HorizontalAlignment

$ws.Range("F18").NumberFormat = $ws.Range("F17").NumberFormat
$ws.Range("F18").HorizontalAlignment = $ws.Range("F17").HorizontalAlignment
$ws.Range("F18").WrapText = $true

$ws.Range("G18").NumberFormat = $ws.Range("G17").NumberFormat
$ws.Range("G18").HorizontalAlignment = $ws.Range("G17").HorizontalAlignment

$ws.Range("H18").NumberFormat = $ws.Range("H17").NumberFormat
$ws.Range("H18").HorizontalAlignment = $ws.Range("H17").HorizontalAlignment
$ws.Range("H18").WrapText = $true
$ws.Range("H18").Font.Bold = $true

$ws.Range("I18").NumberFormat = $ws.Range("I17").NumberFormat
$ws.Range("I18").HorizontalAlignment = $ws.Range("I17").HorizontalAlignment
$ws.Range("I18").WrapText = $true

$ws.Range("K18:N18").NumberFormat = "[$$-409]#,##0.00;[RED]\-[$$-409]#,##0.00"
$ws.Range("K18:N18").HorizontalAlignment = -4108

$ws.Rows.Item(18).RowHeight = 28.35

# --- Column width tweaks for K:M and N ---
$ws.Range("K1:M1").EntireColumn.ColumnWidth = 7.48987854251012
$ws.Range("N1").EntireColumn.ColumnWidth = 8.62753036437247

# --- Sheet view: zoom + selection ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("I15").Select()
